$d = $word.ActiveDocument

# =====================================================================
# Edit 1: the empty paragraph preceding "<<judgeName>>" is removed, and
# the "<<judgeName>>" paragraph (previously styled Heading1) becomes a
# directly-formatted, centered, bold, 14pt (sz=28 half-points) Arial
# paragraph.
# =====================================================================

$rngJudge = $d.Content
$foundJudge = $rngJudge.Find.Execute("<<judgeName>>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($foundJudge) {
    $paras = $d.Paragraphs
    $count = $paras.Count
    $judgeIdx = -1
    for ($i = 1; $i -le $count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Start -le $rngJudge.Start -and $p.Range.End -ge $rngJudge.End) {
            $judgeIdx = $i
            break
        }
    }

    if ($judgeIdx -gt 1) {
        $prevPara = $paras.Item($judgeIdx - 1)
        if ($prevPara.Range.Text -eq "`r") {
            # Delete the empty paragraph entirely (its Range includes the
            # paragraph mark, so this removes the whole paragraph node).
            $prevPara.Range.Delete()
        }
    }

    # Re-resolve the judgeName paragraph after the deletion above.
    $rngJudge2 = $d.Content
    $null = $rngJudge2.Find.Execute("<<judgeName>>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $paras2 = $d.Paragraphs
    $count2 = $paras2.Count
    $judgeIdx2 = -1
    for ($i = 1; $i -le $count2; $i++) {
        $p = $paras2.Item($i)
        if ($p.Range.Start -le $rngJudge2.Start -and $p.Range.End -ge $rngJudge2.End) {
            $judgeIdx2 = $i
            break
        }
    }

    $judgePara = $paras2.Item($judgeIdx2)

    $judgeXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="5" w:line="249" w:lineRule="auto"/><w:ind w:left="118" w:hanging="10"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>&lt;&lt;judgeName&gt;&gt;</w:t></w:r></w:p>
'@

    $judgePara.Range.InsertXML($judgeXml)
}

# =====================================================================
# Edit 2: the "<<es_>>" / "<<cs_{hasPersonalInjury=true}>>" / "Expert
# evidence" paragraph trio loses its direct "es-ES" language tagging,
# and the "E" + "xpert evidence" runs merge into a single run.
# =====================================================================

$rngExpert = $d.Content
$foundExpert = $rngExpert.Find.Execute("hasPersonalInjury=true", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($foundExpert) {
    $paras3 = $d.Paragraphs
    $count3 = $paras3.Count
    $csIdx = -1
    for ($i = 1; $i -le $count3; $i++) {
        $p = $paras3.Item($i)
        if ($p.Range.Start -le $rngExpert.Start -and $p.Range.End -ge $rngExpert.End) {
            $csIdx = $i
            break
        }
    }

    $esIdx = $csIdx - 1
    $expertIdx = $csIdx + 1

    $esPara = $paras3.Item($esIdx)
    $expertPara = $paras3.Item($expertIdx)

    $blockRng = $d.Range($esPara.Range.Start, $expertPara.Range.End)

    $blockXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/></w:rPr><w:t>&lt;&lt;es_&gt;&gt;</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/></w:rPr><w:t>&lt;&lt;cs_{</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>hasPersonalInjury=true</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>}&gt;&gt;</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/></w:rPr><w:t>Expert evidence</w:t></w:r></w:p>
'@

    $blockRng.InsertXML($blockXml)
}

Write-Output "edit complete"
